$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2504.3641875225248
$ws.Range("B1").Value = 1675.4817733253944
$ws.Range("C1").Value = 1659.175700557744
$ws.Range("A2").Value = 2262.5551341635014
$ws.Range("B2").Value = 1502.7606434992442
$ws.Range("C2").Value = 1376.2941818763427
$ws.Range("A3").Value = 2560.3680356691543
$ws.Range("B3").Value = 1740.6138750621997
$ws.Range("C3").Value = 1579.0404426421503
$ws.Range("A4").Value = 2494.2816724362074
$ws.Range("B4").Value = 1886.0880077709817
$ws.Range("C4").Value = 1890.1705214780602
$ws.Range("A5").Value = 2528.8038651354173
$ws.Range("B5").Value = 1752.8912484338457
$ws.Range("C5").Value = 1754.8834586714404
$ws.Range("A6").Value = 2478.7488923483206
$ws.Range("B6").Value = 1826.3368947729757
$ws.Range("C6").Value = 1881.3871662939919
$ws.Range("A7").Value = 2397.2088800618067
$ws.Range("B7").Value = 1848.4308520432855
$ws.Range("C7").Value = 1675.6638407780397
$ws.Range("A8").Value = 2467.996998375711
$ws.Range("B8").Value = 1930.0307207430074
$ws.Range("C8").Value = 1785.2198927447378
$ws.Range("A9").Value = 2643.8788680913858
$ws.Range("B9").Value = 1943.4610312448535
$ws.Range("C9").Value = 1635.4856174790666
$ws.Range("A10").Value = 2394.9108173457316
$ws.Range("B10").Value = 1520.2078910823648
$ws.Range("C10").Value = 1429.4701744687366
$ws.Range("A11").Value = 2163.0433237002462
$ws.Range("B11").Value = 1575.4311100229684
$ws.Range("C11").Value = 1407.2570777772714
$ws.Range("A12").Value = 2770.6264484910107
$ws.Range("B12").Value = 2193.8789571331267
$ws.Range("C12").Value = 1875.0922041911349
$ws.Range("A13").Value = 2528.1018186655519
$ws.Range("B13").Value = 1948.3587415276781
$ws.Range("C13").Value = 1754.633454248107
$ws.Range("A14").Value = 2609.7077573736137
$ws.Range("B14").Value = 2018.5644302498422
$ws.Range("C14").Value = 1759.3304433006706
$ws.Range("A15").Value = 2498.8845153969119
$ws.Range("B15").Value = 2026.4017199517746
$ws.Range("C15").Value = 1838.8582161258253
$ws.Range("A16").Value = 2592.3756000667886
$ws.Range("B16").Value = 1796.2617606532046
$ws.Range("C16").Value = 1554.7405201260492
$ws.Range("A17").Value = 2355.9502413785494
$ws.Range("B17").Value = 1785.7963220480617
$ws.Range("C17").Value = 1699.7042402819582
$ws.Range("A18").Value = 2615.6965196820511
$ws.Range("B18").Value = 2162.0331228077912
$ws.Range("C18").Value = 2083.448066940392
$ws.Range("A19").Value = 2011.2296828283734
$ws.Range("B19").Value = 2033.9312076960384
$ws.Range("C19").Value = 1981.0689856943163
$ws.Range("A20").Value = 2576.2845482917992
$ws.Range("B20").Value = 1981.3402408926472
$ws.Range("C20").Value = 1870.1560410687173
$ws.Range("A21").Value = 2713.0802000001959
$ws.Range("B21").Value = 2020.8570413855018
$ws.Range("C21").Value = 1980.6742389564292
$ws.Range("A22").Value = 2569.5114382281931
$ws.Range("B22").Value = 1982.4267725062407
$ws.Range("C22").Value = 1732.1200200450439
